$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows 18-24: additional aircraft-stability-derivative parameters.
# ---------------------------------------------------------------------------

# Write the (still-plain) label text first. A18 and A22 share identical text
# ("Oswald efficiency factor - (this is an estimation)") - write both while
# plain so the "(this is an estimation)" portion of A18 can be coloured red
# afterwards without disturbing A22's (plain) copy.
$ws.Range("A18").Value = "Oswald efficiency factor - (this is an estimation)"
$ws.Range("A22").Value = "Oswald efficiency factor - (this is an estimation)"

$ws.Range("A19").Value = "Center of gravity of wing"
$ws.Range("B19").Value = "CG_w"

$ws.Range("A20").Value = "Center of gravity of tail"
$ws.Range("B20").Value = "CG_t"

$ws.Range("D18").Value = "<- PLACEHOLDER"

$ws.Range("A21").Value = "Elevator effectiveness ration - (this is an estimation)"
$ws.Range("B21").Value = "T_e"

$ws.Range("B18").Value = "e_w"
$ws.Range("B22").Value = "e_t"

$ws.Range("A23").Value = "Tail span"
$ws.Range("B23").Value = "b_t"

$ws.Range("A24").Value = "Tail chord"
$ws.Range("B24").Value = "c_t"

# Colour the trailing "(this is an estimation)" portion of the two
# estimation labels in red, matching the source workbook's annotation style.
$ws.Range("A18").Characters(28, 23).Font.Color = 255
$ws.Range("A21").Characters(33, 23).Font.Color = 255

# Numeric values for the new rows.
$ws.Range("C18").Value = 0.8
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("C21").Value = 0.665
$ws.Range("C22").Value = 0.8
$ws.Range("C23").Value = 3.04
$ws.Range("C24").Value = 0.833

# "<- PLACEHOLDER" note repeated in column D for rows 19-22.
$ws.Range("D19").Value = "<- PLACEHOLDER"
$ws.Range("D20").Value = "<- PLACEHOLDER"
$ws.Range("D21").Value = "<- PLACEHOLDER"
$ws.Range("D22").Value = "<- PLACEHOLDER"

$ws.Range("E9").Select() | Out-Null
